$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows right before the current row 101, shifting all
# subsequent rows (old 101-181) down to 104-184.
$ws.Range("A101:A103").EntireRow.Insert()

# New row 101: Sandia / Especial, fecha 44596 (2022-02-04)
$ws.Cells.Item(101, "A").Value = 5
$ws.Cells.Item(101, "B").Value = "Macroferia Regional de Talca"
$ws.Cells.Item(101, "C").Value = "Maule"
$ws.Cells.Item(101, "D").Value = 44596
$ws.Cells.Item(101, "E").Value = 7
$ws.Cells.Item(101, "F").Value = 100112028
$ws.Cells.Item(101, "G").Value = "Sandia"
$ws.Cells.Item(101, "H").Value = "Sin especificar"
$ws.Cells.Item(101, "I").Value = "Especial"
$ws.Cells.Item(101, "J").Value = 2000
$ws.Cells.Item(101, "K").Value = 2000
$ws.Cells.Item(101, "L").Value = 2000
$ws.Cells.Item(101, "M").Value = 2000
$ws.Cells.Item(101, "N").Value = "$/unidad"
$ws.Cells.Item(101, "O").Value = "Región del Maule"
$ws.Cells.Item(101, "P").Value = 2000
$ws.Cells.Item(101, "Q").Value = 1
$ws.Cells.Item(101, "R").Value = "Hortaliza"

# New row 102: Sandia / Primera, fecha 44596
$ws.Cells.Item(102, "A").Value = 5
$ws.Cells.Item(102, "B").Value = "Macroferia Regional de Talca"
$ws.Cells.Item(102, "C").Value = "Maule"
$ws.Cells.Item(102, "D").Value = 44596
$ws.Cells.Item(102, "E").Value = 7
$ws.Cells.Item(102, "F").Value = 100112028
$ws.Cells.Item(102, "G").Value = "Sandia"
$ws.Cells.Item(102, "H").Value = "Sin especificar"
$ws.Cells.Item(102, "I").Value = "Primera"
$ws.Cells.Item(102, "J").Value = 3000
$ws.Cells.Item(102, "K").Value = 1500
$ws.Cells.Item(102, "L").Value = 1500
$ws.Cells.Item(102, "M").Value = 1500
$ws.Cells.Item(102, "N").Value = "$/unidad"
$ws.Cells.Item(102, "O").Value = "Región del Maule"
$ws.Cells.Item(102, "P").Value = 1500
$ws.Cells.Item(102, "Q").Value = 1
$ws.Cells.Item(102, "R").Value = "Hortaliza"

# New row 103: Sandia / Segunda, fecha 44596
$ws.Cells.Item(103, "A").Value = 5
$ws.Cells.Item(103, "B").Value = "Macroferia Regional de Talca"
$ws.Cells.Item(103, "C").Value = "Maule"
$ws.Cells.Item(103, "D").Value = 44596
$ws.Cells.Item(103, "E").Value = 7
$ws.Cells.Item(103, "F").Value = 100112028
$ws.Cells.Item(103, "G").Value = "Sandia"
$ws.Cells.Item(103, "H").Value = "Sin especificar"
$ws.Cells.Item(103, "I").Value = "Segunda"
$ws.Cells.Item(103, "J").Value = 3000
$ws.Cells.Item(103, "K").Value = 1000
$ws.Cells.Item(103, "L").Value = 1000
$ws.Cells.Item(103, "M").Value = 1000
$ws.Cells.Item(103, "N").Value = "$/unidad"
$ws.Cells.Item(103, "O").Value = "Región del Maule"
$ws.Cells.Item(103, "P").Value = 1000
$ws.Cells.Item(103, "Q").Value = 1
$ws.Cells.Item(103, "R").Value = "Hortaliza"
